# Auto-generated edits applying the Shinryu_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 199.08696
$ws.Range("I28").Value = 189.19048
$ws.Range("J28").Value = 303
$ws.Range("K28").Value = 189.19048
$ws.Range("L28").Value = 303
$ws.Range("M28").Value = 295.80952
$ws.Range("N28").Value = -1273
$ws.Range("H62").Value = 2096.724
$ws.Range("I62").Value = 2000.2
$ws.Range("K62").Value = 2000.2
$ws.Range("M62").Value = -1376.2
$ws.Range("H65").Value = 2096.724
$ws.Range("I65").Value = 2000.2
$ws.Range("K65").Value = 10001
$ws.Range("M65").Value = -6881
$ws.Range("H98").Value = 769.6286
$ws.Range("I98").Value = 769.6286
$ws.Range("K98").Value = 769.6286
$ws.Range("M98").Value = 728.3714
$ws.Range("H107").Value = 581.6087
$ws.Range("I107").Value = 456.78946
$ws.Range("J107").Value = 1174.5
$ws.Range("K107").Value = 456.78946
$ws.Range("L107").Value = 1174.5
$ws.Range("M107").Value = 1463.21054
$ws.Range("N107").Value = -5014.5
$ws.Range("H111").Value = 200007100
$ws.Range("I111").Value = 250008820
$ws.Range("J111").Value = 250
$ws.Range("K111").Value = 750026460
$ws.Range("L111").Value = 750
$ws.Range("M111").Value = -750023393
$ws.Range("N111").Value = -6884
$ws.Range("H113").Value = 1742.2727
$ws.Range("I113").Value = 1568.9
$ws.Range("J113").Value = 1817.6522
$ws.Range("K113").Value = 1568.9
$ws.Range("L113").Value = 1817.6522
$ws.Range("M113").Value = 1685.1
$ws.Range("N113").Value = -8325.6522
$ws.Range("H115").Value = 10000487
$ws.Range("I115").Value = 10000487
$ws.Range("K115").Value = 30001461
$ws.Range("M115").Value = -29999894
$ws.Range("H118").Value = 418.5
$ws.Range("I118").Value = 418.5
$ws.Range("K118").Value = 1255.5
$ws.Range("M118").Value = 401.5
$ws.Range("H122").Value = 769.6286
$ws.Range("I122").Value = 769.6286
$ws.Range("K122").Value = 2308.8858
$ws.Range("M122").Value = 141.1142
$ws.Range("H127").Value = 1016.8913
$ws.Range("I127").Value = 415.16666
$ws.Range("J127").Value = 1107.15
$ws.Range("K127").Value = 1245.49998
$ws.Range("L127").Value = 3321.45
$ws.Range("M127").Value = 3714.50002
$ws.Range("N127").Value = -13241.45
$ws.Range("H129").Value = 1006.42
$ws.Range("I129").Value = 390.41666
$ws.Range("J129").Value = 1090.4204
$ws.Range("K129").Value = 1171.24998
$ws.Range("L129").Value = 3271.2612
$ws.Range("M129").Value = 3828.75002
$ws.Range("N129").Value = -13271.2612
$ws.Range("H131").Value = 27269.975
$ws.Range("I131").Value = 31430
$ws.Range("J131").Value = 4389.8335
$ws.Range("K131").Value = 94290
$ws.Range("L131").Value = 13169.5005
$ws.Range("M131").Value = -89250
$ws.Range("N131").Value = -23249.5005
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3730.625
$ws.Range("I2").Value = 4238.6665
$ws.Range("J2").Value = 2206.5
$ws.Range("K2").Value = 4238.6665
$ws.Range("L2").Value = 2206.5
$ws.Range("M2").Value = -4125.6665
$ws.Range("N2").Value = -2432.5
$ws.Range("H45").Value = 2205.1667
$ws.Range("I45").Value = 1914.8889
$ws.Range("J45").Value = 3076
$ws.Range("K45").Value = 1914.8889
$ws.Range("L45").Value = 3076
$ws.Range("M45").Value = -1537.8889
$ws.Range("N45").Value = -3830
$ws.Range("H110").Value = 990.35486
$ws.Range("I110").Value = 788.6818
$ws.Range("J110").Value = 1483.3334
$ws.Range("K110").Value = 788.6818
$ws.Range("L110").Value = 1483.3334
$ws.Range("M110").Value = 1256.3182
$ws.Range("N110").Value = -5573.3334
$ws.Range("H116").Value = 3730.625
$ws.Range("I116").Value = 4238.6665
$ws.Range("J116").Value = 2206.5
$ws.Range("K116").Value = 4238.6665
$ws.Range("L116").Value = 2206.5
$ws.Range("M116").Value = -1944.6665
$ws.Range("N116").Value = -6794.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3730.625
$ws.Range("I3").Value = 4238.6665
$ws.Range("J3").Value = 2206.5
$ws.Range("K3").Value = 4238.6665
$ws.Range("L3").Value = 2206.5
$ws.Range("M3").Value = -4124.6665
$ws.Range("N3").Value = -2434.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2772.8
$ws.Range("I16").Value = 1601.5
$ws.Range("J16").Value = 4111.4287
$ws.Range("K16").Value = 1601.5
$ws.Range("L16").Value = 4111.4287
$ws.Range("M16").Value = -1314.5
$ws.Range("N16").Value = -4685.4287
$ws.Range("H22").Value = 372.58334
$ws.Range("I22").Value = 221.125
$ws.Range("J22").Value = 675.5
$ws.Range("K22").Value = 221.125
$ws.Range("L22").Value = 675.5
$ws.Range("M22").Value = 128.875
$ws.Range("N22").Value = -1375.5
$ws.Range("H107").Value = 527.7568
$ws.Range("I107").Value = 418.17856
$ws.Range("K107").Value = 418.17856
$ws.Range("M107").Value = 1501.82144
$ws.Range("H113").Value = 2772.8
$ws.Range("I113").Value = 1601.5
$ws.Range("J113").Value = 4111.4287
$ws.Range("K113").Value = 1601.5
$ws.Range("L113").Value = 4111.4287
$ws.Range("M113").Value = 568.5
$ws.Range("N113").Value = -8451.4287
$ws.Range("H122").Value = 758.1852
$ws.Range("I122").Value = 804.7273
$ws.Range("J122").Value = 553.4
$ws.Range("K122").Value = 2414.1819
$ws.Range("L122").Value = 1660.2
$ws.Range("M122").Value = 35.81809999999996
$ws.Range("N122").Value = -6560.2
$ws.Range("H132").Value = 2234.7144
$ws.Range("I132").Value = 1580.7826
$ws.Range("J132").Value = 5242.8
$ws.Range("K132").Value = 4742.3478
$ws.Range("L132").Value = 15728.4
$ws.Range("M132").Value = -2212.3478
$ws.Range("N132").Value = -20788.4
$ws.Range("H134").Value = 2352.3572
$ws.Range("I134").Value = 1446.7273
$ws.Range("J134").Value = 5673
$ws.Range("K134").Value = 4340.1819
$ws.Range("L134").Value = 17019
$ws.Range("M134").Value = -1805.1819
$ws.Range("N134").Value = -22089
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 500
$ws.Range("I23").Value = 400
$ws.Range("J23").Value = 600
$ws.Range("K23").Value = 1200
$ws.Range("L23").Value = 1800
$ws.Range("M23").Value = -965
$ws.Range("N23").Value = -2270
$ws.Range("H74").Value = 5015
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 5015
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 15045
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -17167
$ws.Range("H77").Value = 5015
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 5015
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 45135
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -55743
$ws.Range("H131").Value = 859.1111
$ws.Range("J131").Value = 943.1429000000001
$ws.Range("L131").Value = 2829.4287
$ws.Range("N131").Value = -12909.4287
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1708.92
$ws.Range("I102").Value = 1640.1305
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1640.1305
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -18.13049999999998
$ws.Range("N102").Value = -5744
$ws.Range("H107").Value = 241.3125
$ws.Range("I107").Value = 241.3125
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 241.3125
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1678.6875
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 17030.285
$ws.Range("I113").Value = 1750
$ws.Range("K113").Value = 1750
$ws.Range("M113").Value = 420
$ws.Range("H122").Value = 12501334
$ws.Range("I122").Value = 16667910
$ws.Range("J122").Value = 1604
$ws.Range("K122").Value = 50003730
$ws.Range("L122").Value = 4812
$ws.Range("M122").Value = -50001280
$ws.Range("N122").Value = -9712
$ws.Range("H126").Value = 3165.3823
$ws.Range("I126").Value = 2837.182
$ws.Range("J126").Value = 3767.0833
$ws.Range("K126").Value = 8511.545999999998
$ws.Range("L126").Value = 11301.2499
$ws.Range("M126").Value = -6041.545999999998
$ws.Range("N126").Value = -16241.2499
$ws.Range("H132").Value = 3747.64
$ws.Range("I132").Value = 3486.9546
$ws.Range("K132").Value = 10460.8638
$ws.Range("M132").Value = -7930.863799999999
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2231.5
$ws.Range("I7").Value = 2212.375
$ws.Range("J7").Value = 2250.625
$ws.Range("K7").Value = 2212.375
$ws.Range("L7").Value = 2250.625
$ws.Range("M7").Value = -2100.375
$ws.Range("N7").Value = -2474.625
$ws.Range("H94").Value = 31107.666
$ws.Range("J94").Value = 31107.666
$ws.Range("L94").Value = 31107.666
$ws.Range("N94").Value = -32459.666
$ws.Range("H126").Value = 2231.5
$ws.Range("I126").Value = 2212.375
$ws.Range("J126").Value = 2250.625
$ws.Range("K126").Value = 6637.125
$ws.Range("L126").Value = 2250.625
$ws.Range("M126").Value = -4167.125
$ws.Range("N126").Value = -11691.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 556
$ws.Range("I107").Value = 495
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1485
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = 435
$ws.Range("N107").Value = -6240
$ws.Range("H113").Value = 874.8333
$ws.Range("I113").Value = 861.875
$ws.Range("J113").Value = 900.75
$ws.Range("K113").Value = 2585.625
$ws.Range("L113").Value = 2702.25
$ws.Range("M113").Value = -415.625
$ws.Range("N113").Value = -7042.25
$ws.Range("H122").Value = 2253.0645
$ws.Range("I122").Value = 2095.8333
$ws.Range("J122").Value = 2470.7693
$ws.Range("K122").Value = 6287.499899999999
$ws.Range("L122").Value = 7412.3079
$ws.Range("M122").Value = -3837.499899999999
$ws.Range("N122").Value = -12312.3079
$ws.Range("H126").Value = 1894.08
$ws.Range("I126").Value = 1968.9048
$ws.Range("J126").Value = 1501.25
$ws.Range("K126").Value = 5906.7144
$ws.Range("L126").Value = 1501.25
$ws.Range("M126").Value = -3436.7144
$ws.Range("N126").Value = -9443.75
